# Apply the acs5_variable_sheets agexnum_disabilities_poverty update:
#   - rename the three sheets to append " - done"
#   - make "Age by disability by poverty" the active sheet/tab
#   - update cell selections on sheet 1 and sheet 2
#   - refresh the wrapped-text row heights on all three sheets

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Age by number of disabilities
$ws2 = $wb.Worksheets.Item(2)   # Age by disability by poverty
$ws3 = $wb.Worksheets.Item(3)   # Age EJS

# --- Sheet 1: selection + row heights -------------------------------------
[void]$ws1.Range("C15:C16").Select()

$ws1Heights = @{
    1 = 34
    2 = 29
    3 = 29
    4 = 43
    5 = 51
    6 = 29
    7 = 29
    8 = 43
    9 = 51
    10 = 29
    11 = 29
    12 = 43
    13 = 51
    14 = 29
}
foreach ($r in $ws1Heights.Keys) {
    $ws1.Rows.Item($r).RowHeight = $ws1Heights[$r]
}

# --- Sheet 2: row heights ---------------------------------------------------
$ws2Heights = @{
    1 = 51
    2 = 43
    3 = 43
    4 = 68
    5 = 71
    6 = 71
    7 = 43
    8 = 57
    9 = 57
    10 = 43
    11 = 68
    12 = 71
    13 = 71
    14 = 43
    15 = 57
    16 = 57
    17 = 43
    18 = 68
    19 = 68
    20 = 68
    21 = 43
    22 = 57
    23 = 57
}
foreach ($r in $ws2Heights.Keys) {
    $ws2.Rows.Item($r).RowHeight = $ws2Heights[$r]
}

# --- Sheet 3: row heights ----------------------------------------------------
$ws3Heights = @{
    1 = 18
    2 = 41
    3 = 21
    4 = 21
    5 = 21
    6 = 21
}
foreach ($r in $ws3Heights.Keys) {
    $ws3.Rows.Item($r).RowHeight = $ws3Heights[$r]
}

# --- Sheet 2: become the active tab + selection -----------------------------
# (done after the row-height/selection edits above so sheet 2 ends up the
# active/tabSelected sheet, matching the target workbook view)
[void]$ws2.Activate()
[void]$ws2.Range("I11").Select()

# --- Rename the sheets (append " - done") -----------------------------------
$ws1.Name = "Age by num disabilities - done"
$ws2.Name = "Age x disabil x poverty - done"
$ws3.Name = "Age EJS - done"
